$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) to reflect the new commit.
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection left behind by the edit.
$ws.Range("E8").Select()
